$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete rows 19-33 (other workers), shifting the signature block up from 38/39 to 23/24
$ws.Range("A19:A33").EntireRow.Delete()

# Update the 3 remaining worker rows (16-18) with Cesar Enrique Polo Marimon's data
$ws.Range("C16").Value = "1047396061"
$ws.Range("D16").Value = "CESAR ENRIQUE POLO MARIMON"
$ws.Range("E16").Value = "1904"
$ws.Range("F16").Value = 31249
$ws.Range("G16").Value = 908526

$ws.Range("C17").Value = "1047396061"
$ws.Range("D17").Value = "CESAR ENRIQUE POLO MARIMON"
$ws.Range("E17").Value = "1905"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 908526

$ws.Range("C18").Value = "1047396061"
$ws.Range("D18").Value = "CESAR ENRIQUE POLO MARIMON"
$ws.Range("E18").Value = "1906"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 908526

# Update summary fields
$ws.Range("E11").Value = 93747
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 3
